$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68:I70").Copy()
$ws.Range("A96").PasteSpecial(-4122)

$ws.Range("J84").Copy()
$ws.Range("J96").PasteSpecial(-4122)
$ws.Range("J87").Copy()
$ws.Range("J97").PasteSpecial(-4122)
$ws.Range("J85").Copy()
$ws.Range("J98").PasteSpecial(-4122)

$rows = 96,97,98
$times = "2024-10-01T19:27:47","2024-10-05T00:06:07","2024-10-11T01:31:06"
for ($i = 0; $i -lt 3; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value2 = 85
    $ws.Cells.Item($r, 2).Value2 = "Phobos"
    $ws.Cells.Item($r, 3).Value2 = "LNO centre"
    $ws.Cells.Item($r, 4).Value2 = 60
    $ws.Cells.Item($r, 5).Value2 = 6
    $ws.Cells.Item($r, 6).Value2 = "CARBONATES 174 175 176 189 190 191 500MS"
    $ws.Cells.Item($r, 7).Value2 = 12
    $ws.Cells.Item($r, 8).Value2 = 3922
    $ws.Cells.Item($r, 9).Value2 = 235
    $ws.Cells.Item($r, 10).Value2 = $times[$i]
}
$ws.Range("F89").Select()
Write-Host "done"
